{"js": "// Update the worksheet date and all twenty-five \"three-digit \u00d7 one-digit\"\n// multiplication prompts to the next day's freshly generated values.\nconst replacements = [\n  [\"2024-12-25 Wednesday\", \"2024-12-26 Thursday\"],\n  [\"554\u00d72=\", \"677\u00d76=\"],\n  [\"797\u00d76=\", \"257\u00d73=\"],\n  [\"861\u00d77=\", \"603\u00d76=\"],\n  [\"862\u00d74=\", \"201\u00d73=\"],\n  [\"963\u00d73=\", \"463\u00d78=\"],\n  [\"801\u00d76=\", \"254\u00d79=\"],\n  [\"800\u00d75=\", \"453\u00d79=\"],\n  [\"604\u00d76=\", \"631\u00d74=\"],\n  [\"720\u00d74=\", \"521\u00d75=\"],\n  [\"906\u00d74=\", \"825\u00d75=\"],\n  [\"760\u00d73=\", \"440\u00d76=\"],\n  [\"630\u00d72=\", \"357\u00d75=\"],\n  [\"270\u00d74=\", \"869\u00d75=\"],\n  [\"295\u00d75=\", \"937\u00d73=\"],\n  [\"570\u00d79=\", \"102\u00d78=\"],\n  [\"595\u00d72=\", \"645\u00d77=\"],\n  [\"893\u00d73=\", \"667\u00d73=\"],\n  [\"753\u00d76=\", \"180\u00d77=\"],\n  [\"492\u00d73=\", \"599\u00d79=\"],\n  [\"562\u00d72=\", \"826\u00d74=\"],\n  [\"121\u00d75=\", \"653\u00d72=\"],\n  [\"129\u00d78=\", \"773\u00d75=\"],\n  [\"568\u00d77=\", \"523\u00d78=\"],\n  [\"636\u00d79=\", \"291\u00d72=\"],\n  [\"450\u00d74=\", \"780\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all twenty-five \"three-digit \u00d7 one-digit\"\n# multiplication prompts to the next day's freshly generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-25 Wednesday\", \"2024-12-26 Thursday\"),\n    @(\"554\u00d72=\", \"677\u00d76=\"),\n    @(\"797\u00d76=\", \"257\u00d73=\"),\n    @(\"861\u00d77=\", \"603\u00d76=\"),\n    @(\"862\u00d74=\", \"201\u00d73=\"),\n    @(\"963\u00d73=\", \"463\u00d78=\"),\n    @(\"801\u00d76=\", \"254\u00d79=\"),\n    @(\"800\u00d75=\", \"453\u00d79=\"),\n    @(\"604\u00d76=\", \"631\u00d74=\"),\n    @(\"720\u00d74=\", \"521\u00d75=\"),\n    @(\"906\u00d74=\", \"825\u00d75=\"),\n    @(\"760\u00d73=\", \"440\u00d76=\"),\n    @(\"630\u00d72=\", \"357\u00d75=\"),\n    @(\"270\u00d74=\", \"869\u00d75=\"),\n    @(\"295\u00d75=\", \"937\u00d73=\"),\n    @(\"570\u00d79=\", \"102\u00d78=\"),\n    @(\"595\u00d72=\", \"645\u00d77=\"),\n    @(\"893\u00d73=\", \"667\u00d73=\"),\n    @(\"753\u00d76=\", \"180\u00d77=\"),\n    @(\"492\u00d73=\", \"599\u00d79=\"),\n    @(\"562\u00d72=\", \"826\u00d74=\"),\n    @(\"121\u00d75=\", \"653\u00d72=\"),\n    @(\"129\u00d78=\", \"773\u00d75=\"),\n    @(\"568\u00d77=\", \"523\u00d78=\"),\n    @(\"636\u00d79=\", \"291\u00d72=\"),\n    @(\"450\u00d74=\", \"780\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
